# Applies the "Actualizacion del log y el task" edit to the LOGT worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 10 for the new log entry, pushing the old row 10 to row 11 ---
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).RowHeight = 30

$ws.Range("A10").Value = 41912
$ws.Range("B10").Value = 0.96875
$ws.Range("C10").Value = 0.99305555555555547
$ws.Range("D10").Value = 5
$ws.Range("E10").Formula = "=((HOUR(C10)-HOUR(B10))*60)+(MINUTE(C10)-MINUTE(B10))-D10"
$ws.Range("H10").Value = "Reunión de equipo para discutir estados de las tareas del ciclo #1."

# --- Update existing Phase/Task (column F) values for rows 6-9 ---
# Row 6: numeric 1 -> text "-" (centered style)
$ws.Range("F6").Value = "-"
$ws.Range("F6").HorizontalAlignment = -4108   # xlCenter

# New row 10 also gets the "-" phase/task marker
$ws.Range("F10").Value = "-"
$ws.Range("F10").HorizontalAlignment = -4108   # xlCenter

# Row 7: 2 -> 1
$ws.Range("F7").Value = 1

# Row 8: 3 -> 7
$ws.Range("F8").Value = 7

# Row 9: 4 -> 7
$ws.Range("F9").Value = 7

# --- Update the row that shifted down from 10 to 11 ---
$ws.Range("A11").Value = 41913
$ws.Range("F11").Value = 8

# --- Update worksheet selection ---
$ws.Range("F21").Select()
